$d = $word.ActiveDocument

# 1. Fix scenario title: "Scenario registrovanja korisnika" -> "Scenario pregled početne strane"
$d.Content.Find.Execute("Scenario registrovanja korisnika", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Scenario pregled početne strane", 2)

# 2. Fix typo: "Alternatini" -> "Alternativni"
$d.Content.Find.Execute("Alternatini", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Alternativni", 2)

# 3. Add missing space after comma: "lajkove,slike" -> "lajkove, slike"
$d.Content.Find.Execute("lajkove,slike", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "lajkove, slike", 2)
